$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1 ("We chose to utilize more of the GitHub features...")
# The original paragraph is split across 3 runs; replacing text that spans
# all three runs merges them back into a single run (taking the formatting
# of the first run), matching the target structure.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "We chose to utilize more of the GitHub features, such as the issues feature and the projects feature that allows you to specify the progress of each task (kind of like a " + [char]0x2018 + "to-do' list), to learn how to use GitHub in a more professional manner and to keep each other informed on our progress of the tasks we were assigned. Every day, or as often as we could, we would update these features. These came in handy when we were unsure of what tasks were needed to be complete or were completed. They also helped us to know who was working on which task and who completed which task. After each scrum meeting, we would go onto GitHub and update these when we could.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "We chose to utilize more of the GitHub features, such as the issues and the projects that allows you to specify the progress of each task (" + [char]0x2018 + "to-do' list), to learn how to use GitHub in a more professional manner and to keep each other informed on our progress of the tasks we were assigned. At least once a day we would update the Kanban board on GitHub and add more issues when relevant. These came in handy when we were unsure of what tasks needed to be complete or were completed. They also helped us to know who was working on which task and who completed which task. After each scrum meeting, the projects and issues where updated with the most up to date information.", `
    2)

# Find.Execute's replacement text runs through smart-quote AutoCorrect, which
# turns the trailing straight apostrophe of "to-do'" into a curly closing
# quote. The source keeps it straight (‘to-do'), so patch that one character
# back via a direct Range.Text assignment, which bypasses AutoCorrect.
$p1Range = $d.Paragraphs(1).Range
$quoteFix = $d.Range($p1Range.Start, $p1Range.End)
$quoteFix.Find.Execute("to-do" + [char]0x2019, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$aposRange = $d.Range($quoteFix.End - 1, $quoteFix.End)
$aposRange.Text = [char]0x27

# ---------------------------------------------------------------------------
# Paragraph 2 ("We also chose to update our sprint backlog every day...")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "sprint backlog correct and not have to guess how long it took us to complete a task.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "sprint backlog correct and more efficient so that we don" + [char]0x2019 + "t have to guess how long it took us to complete a task.", `
    2)

# ---------------------------------------------------------------------------
# Paragraph 3 ("We did not scrum with our backlog up every time...")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "see each other when he/she was speaking.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "see each other when they were speaking.", `
    2)

$d.Content.Find.Execute( `
    "It also allows each person to receive help if they are struggling with a certain task.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This also allowed each person to receive help if they were struggling on a certain task.", `
    2)

# ---------------------------------------------------------------------------
# Paragraph 4 ("We also improved our communication using the issues feature...")
# First rewrite the sentence text (including the trailing new sentence), then
# split the tail into three runs: "...not as ", "individuals", "." -- each
# carrying the same run formatting as the paragraph's run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "We also improved our communication using the issues feature in GitHub and by speaking to each other when we needed help. This allowed for everyone to learn more when they were coding and for the project to be completed in a more timely manner.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "We also improved our communication using the issues feature in GitHub and by speaking to each other when we needed help. This allowed for everyone to learn more when they were coding and for the project to be completed in a timelier manner. This also helped us work better as team and not as individuals.", `
    2)

$p4 = $d.Paragraphs(4)
$p4Range = $d.Range($p4.Range.Start, $p4.Range.End)
$p4Range.Find.Execute("individuals", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$individualsRange = $d.Range($p4Range.Start, $p4Range.End)
# Toggle a format off/on (net no-op) purely to force a run boundary around
# "individuals" so it becomes its own <w:r>.
$individualsRange.Bold = 1
$individualsRange.Bold = 0

$periodRange = $d.Range($individualsRange.End, $individualsRange.End + 1)
# Same trick for the trailing "." so it becomes its own <w:r> too.
$periodRange.Bold = 1
$periodRange.Bold = 0
